$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename property code from $NAME to $$NAME
$ws.Range("A5").Value = '$$NAME'

# Remove "Internal Assignment" column (O) contents from export
$ws.Range("O4:O7").ClearContents()

# Match the new selection recorded for the sheet
$ws.Range("O4:O7").Select()
